# Update NATMI LR-pair TPM-derived metrics for Slit1-Robo2 sheet
# per the "update scripts wuth new tpm" commit: new TPM values propagate
# through ligand/receptor expression, specificity, and edge-weight columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 0.4218206666666666
    "H2" = 1.265462
    "I2" = 0.204479520571209
    "J2" = 0.204479520571209
    "M2" = 1.302860333333333
    "N2" = 3.908581
    "O2" = 0.9669439908960468
    "P2" = 0.9669439908960467
    "Q2" = 0.5495734143802221
    "R2" = 4.946160729421999
    "S2" = 0.1977202436776351
    "T2" = 0.1977202436776351
    "G3" = 0.4218206666666666
    "H3" = 1.265462
    "I3" = 0.204479520571209
    "J3" = 0.204479520571209
    "O3" = 0.008324674682103805
    "P3" = 0.008324674682103805
    "Q3" = 0.00473142181111111
    "R3" = 0.0425827963
    "S3" = 0.001702225487907868
    "T3" = 0.001702225487907867
    "G4" = 0.4218206666666666
    "H4" = 1.265462
    "I4" = 0.204479520571209
    "J4" = 0.204479520571209
    "M4" = 0.033323
    "N4" = 0.099969
    "O4" = 0.02473133442184949
    "P4" = 0.02473133442184949
    "Q4" = 0.01405633007533333
    "R4" = 0.126506970678
    "S4" = 0.005057051405666022
    "T4" = 0.005057051405666021
    "I5" = 0.3030684321645684
    "J5" = 0.3030684321645683
    "M5" = 1.302860333333333
    "N5" = 3.908581
    "O5" = 0.9669439908960468
    "P5" = 0.9669439908960467
    "Q5" = 0.8145478461132221
    "R5" = 7.330930615018999
    "S5" = 0.2930501993118155
    "T5" = 0.2930501993118155
    "I6" = 0.3030684321645684
    "J6" = 0.3030684321645683
    "O6" = 0.008324674682103805
    "P6" = 0.008324674682103805
    "S6" = 0.002522946104185277
    "T6" = 0.002522946104185276
    "I7" = 0.3030684321645684
    "J7" = 0.3030684321645683
    "M7" = 0.033323
    "N7" = 0.099969
    "O7" = 0.02473133442184949
    "P7" = 0.02473133442184949
    "Q7" = 0.02083352849233333
    "R7" = 0.187501756431
    "S7" = 0.007495286748567546
    "T7" = 0.007495286748567544
    "E8" = 3
    "F8" = 1
    "G8" = 0.469433
    "H8" = 1.408299
    "I8" = 0.2275598195290835
    "J8" = 0.2275598195290835
    "M8" = 1.302860333333333
    "N8" = 3.908581
    "O8" = 0.9669439908960468
    "P8" = 0.9669439908960467
    "Q8" = 0.6116056348576666
    "R8" = 5.504450713719
    "S8" = 0.2200376000630362
    "T8" = 0.2200376000630361
    "E9" = 3
    "F9" = 1
    "G9" = 0.469433
    "H9" = 1.408299
    "I9" = 0.2275598195290835
    "J9" = 0.2275598195290835
    "O9" = 0.008324674682103805
    "P9" = 0.008324674682103805
    "Q9" = 0.005265473483333333
    "R9" = 0.04738926135
    "S9" = 0.001894361468297872
    "T9" = 0.001894361468297872
    "E10" = 3
    "F10" = 1
    "G10" = 0.469433
    "H10" = 1.408299
    "I10" = 0.2275598195290835
    "J10" = 0.2275598195290835
    "M10" = 0.033323
    "N10" = 0.099969
    "O10" = 0.02473133442184949
    "P10" = 0.02473133442184949
    "Q10" = 0.015642915859
    "R10" = 0.140786242731
    "S10" = 0.00562785799774948
    "T10" = 0.005627857997749479
    "G11" = 0.546446
    "H11" = 1.639338
    "I11" = 0.2648922277351392
    "J11" = 0.2648922277351391
    "M11" = 1.302860333333333
    "N11" = 3.908581
    "O11" = 0.9669439908960468
    "P11" = 0.9669439908960467
    "Q11" = 0.7119428177086666
    "R11" = 6.407485359378
    "S11" = 0.2561359478435599
    "T11" = 0.2561359478435599
    "G12" = 0.546446
    "H12" = 1.639338
    "I12" = 0.2648922277351392
    "J12" = 0.2648922277351391
    "O12" = 0.008324674682103805
    "P12" = 0.008324674682103805
    "Q12" = 0.006129302633333333
    "R12" = 0.0551637237
    "S12" = 0.002205141621712788
    "T12" = 0.002205141621712788
    "G13" = 0.546446
    "H13" = 1.639338
    "I13" = 0.2648922277351392
    "J13" = 0.2648922277351391
    "M13" = 0.033323
    "N13" = 0.099969
    "O13" = 0.02473133442184949
    "P13" = 0.02473133442184949
    "Q13" = 0.018209220058
    "R13" = 0.163882980522
    "S13" = 0.006551138269866441
    "T13" = 0.006551138269866439
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells with new TPM-derived values"
